# Apply updated cryptocurrency price/volume data per the commit diff.
# Each target cell is forced to the Text number format before assignment,
# and the style is reset back to Normal afterward, so that numeric-looking
# strings (e.g. "313.14", "0.0777") are preserved verbatim as text, exactly
# like the original inlineStr cells, instead of being auto-converted to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '41.486.56'
Set-TextValue 'E2' '  -2.73%  '
Set-TextValue 'D3' '2.469.39'
Set-TextValue 'E3' '  -2.36%  '
Set-TextValue 'E4' '  +0.84%  '
Set-TextValue 'D5' '313.14'
Set-TextValue 'E5' '  -0.43%  '
Set-TextValue 'D6' '92.05'
Set-TextValue 'E6' '  -6.95%  '
Set-TextValue 'D7' '0.543'
Set-TextValue 'E7' '  -3.57%  '
Set-TextValue 'E8' '  +0.81%  '
Set-TextValue 'D9' '0.491'
Set-TextValue 'E9' '  -5.00%  '
Set-TextValue 'D10' '32.76'
Set-TextValue 'E10' '  -7.05%  '
Set-TextValue 'D11' '0.0777'
Set-TextValue 'E11' '  -2.95%  '
Set-TextValue 'E12' '  -0.22%  '
Set-TextValue 'D13' '2.855.91'
Set-TextValue 'E13' '  -2.08%  '
Set-TextValue 'D14' '6.83'
Set-TextValue 'E14' '  -5.34%  '
Set-TextValue 'B15' 'WrappedEther'
Set-TextValue 'C15' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D15' '2.555.74'
Set-TextValue 'E15' '  +1.17%  '
Set-TextValue 'B16' 'Chainlink'
Set-TextValue 'C16' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D16' '15.22'
Set-TextValue 'E16' '  -0.15%  '
Set-TextValue 'D17' '0.779'
Set-TextValue 'E17' '  -3.91%  '
Set-TextValue 'D18' '41.345.35'
Set-TextValue 'E18' '  -3.08%  '
Set-TextValue 'D19' '6.26'
Set-TextValue 'E19' '  -5.29%  '
Set-TextValue 'D20' '0.0₃0916'
Set-TextValue 'E20' '  -2.55%  '
Set-TextValue 'D21' '69.81'
Set-TextValue 'E21' '  +1.08%  '
Set-TextValue 'D22' '10.96'
Set-TextValue 'E22' '  -10.13%  '
Set-TextValue 'D23' '234.10'
Set-TextValue 'E23' '  -3.32%  '
Set-TextValue 'D24' '2.73'
Set-TextValue 'E24' '  -4.38%  '
Set-TextValue 'E25' '  -0.17%  '
Set-TextValue 'D26' '1.86'
Set-TextValue 'E26' '  -6.50%  '
Set-TextValue 'D27' '23.86'
Set-TextValue 'E27' '  -6.35%  '
Set-TextValue 'D28' '2.25'
Set-TextValue 'E28' '  -0.30%  '
Set-TextValue 'D29' '9.68'
Set-TextValue 'E29' '  -3.16%  '
Set-TextValue 'D30' '36.03'
Set-TextValue 'E30' '  -4.58%  '
Set-TextValue 'D31' '152.28'
Set-TextValue 'E31' '  -2.42%  '
Set-TextValue 'D32' '5.40'
Set-TextValue 'E32' '  -8.87%  '
Set-TextValue 'D33' '2.56'
Set-TextValue 'E33' '  -3.01%  '
Set-TextValue 'D34' '2.55'
Set-TextValue 'E34' '  -5.74%  '
Set-TextValue 'D35' '0.0746'
Set-TextValue 'E35' '  -4.79%  '
Set-TextValue 'D36' '17.52'
Set-TextValue 'E36' '  -0.32%  '
Set-TextValue 'D37' '2.98'
Set-TextValue 'E37' '  -4.77%  '
Set-TextValue 'D38' '1.84'
Set-TextValue 'E38' '  -6.73%  '
Set-TextValue 'D39' '0.112'
Set-TextValue 'E39' '  -3.90%  '
Set-TextValue 'D40' '0.0997'
Set-TextValue 'E40' '  -7.90%  '
Set-TextValue 'D41' '4.01'
Set-TextValue 'E41' '  -5.25%  '
Set-TextValue 'B42' 'FirstDigitalUSD'
Set-TextValue 'C42' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D42' '1.01'
Set-TextValue 'E42' '  +1.01%  '
Set-TextValue 'B43' 'EnergySwap'
Set-TextValue 'C43' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D43' '19.27'
Set-TextValue 'E43' '  -12.22%  '
Set-TextValue 'D44' '1.959.33'
Set-TextValue 'E44' '  -3.39%  '
Set-TextValue 'E45' '  -5.20%  '
Set-TextValue 'D46' '2.93'
Set-TextValue 'E46' '  -8.94%  '
Set-TextValue 'D47' '8.70'
Set-TextValue 'E47' '  -2.56%  '
Set-TextValue 'D48' '2.719.82'
Set-TextValue 'E48' '  -1.79%  '
Set-TextValue 'D49' '68.08'
Set-TextValue 'E49' '  -5.43%  '
Set-TextValue 'D50' '95.19'
Set-TextValue 'E50' '  -5.04%  '
Set-TextValue 'D51' '0.175'
Set-TextValue 'E51' '  -7.23%  '
